$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4313.478
$ws.Range("I62").Value = 5362.3076
$ws.Range("J62").Value = 2950
$ws.Range("K62").Value = 5362.3076
$ws.Range("L62").Value = 2950
$ws.Range("M62").Value = -4738.3076
$ws.Range("N62").Value = -4198
$ws.Range("H65").Value = 4313.478
$ws.Range("I65").Value = 5362.3076
$ws.Range("J65").Value = 2950
$ws.Range("K65").Value = 26811.538
$ws.Range("L65").Value = 14750
$ws.Range("M65").Value = -23691.538
$ws.Range("N65").Value = -20990
$ws.Range("H107").Value = 1350
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 1700
$ws.Range("K107").Value = 1000
$ws.Range("L107").Value = 1700
$ws.Range("M107").Value = 920
$ws.Range("N107").Value = -5540
$ws.Range("H132").Value = 18684.611
$ws.Range("I132").Value = 2656.375
$ws.Range("K132").Value = 7969.125
$ws.Range("M132").Value = -5439.125
$ws.Range("H137").Value = 1482267.9
$ws.Range("I137").Value = 2263268.2
$ws.Range("K137").Value = 6789804.600000001
$ws.Range("M137").Value = -6787254.600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13525.191
$ws.Range("I32").Value = 12867.25
$ws.Range("K32").Value = 12867.25
$ws.Range("M32").Value = -12580.25
$ws.Range("H122").Value = 1628.826
$ws.Range("I122").Value = 1616.2778
$ws.Range("J122").Value = 1674
$ws.Range("K122").Value = 4848.8334
$ws.Range("L122").Value = 5022
$ws.Range("M122").Value = -2398.8334
$ws.Range("N122").Value = -9922
$ws.Range("H132").Value = 14287601
$ws.Range("I132").Value = 18519742
$ws.Range("J132").Value = 4124.5
$ws.Range("K132").Value = 55559226
$ws.Range("L132").Value = 12373.5
$ws.Range("M132").Value = -55556696
$ws.Range("N132").Value = -17433.5
$ws.Range("H133").Value = 42933.332
$ws.Range("J133").Value = 42933.332
$ws.Range("L133").Value = 42933.332
$ws.Range("N133").Value = -47993.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2700
$ws.Range("I86").Value = 3000
$ws.Range("J86").Value = 2500
$ws.Range("K86").Value = 3000
$ws.Range("L86").Value = 2500
$ws.Range("M86").Value = -1877
$ws.Range("N86").Value = -4746
$ws.Range("H89").Value = 2700
$ws.Range("I89").Value = 3000
$ws.Range("J89").Value = 2500
$ws.Range("K89").Value = 15000
$ws.Range("L89").Value = 12500
$ws.Range("M89").Value = -9384
$ws.Range("N89").Value = -23732
$ws.Range("H105").Value = 2332.5
$ws.Range("I105").Value = 1767.35
$ws.Range("J105").Value = 3139.8572
$ws.Range("K105").Value = 1767.35
$ws.Range("L105").Value = 3139.8572
$ws.Range("M105").Value = -20.34999999999991
$ws.Range("N105").Value = -6633.8572
$ws.Range("H118").Value = 20000
$ws.Range("J118").Value = 20000
$ws.Range("L118").Value = 20000
$ws.Range("N118").Value = -23314

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 4168.7
$ws.Range("I86").Value = 4168.7
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 4168.7
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -3045.7
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 4168.7
$ws.Range("I89").Value = 4168.7
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 20843.5
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -15227.5
$ws.Range("N89").ClearContents()
$ws.Range("H105").Value = 2414.111
$ws.Range("I105").Value = 2173.3845
$ws.Range("J105").Value = 3040
$ws.Range("K105").Value = 2173.3845
$ws.Range("L105").Value = 3040
$ws.Range("M105").Value = -426.3845000000001
$ws.Range("N105").Value = -6534
$ws.Range("H107").Value = 781.2941
$ws.Range("I107").Value = 513.52
$ws.Range("J107").Value = 1525.1111
$ws.Range("K107").Value = 513.52
$ws.Range("L107").Value = 1525.1111
$ws.Range("M107").Value = 1406.48
$ws.Range("N107").Value = -5365.1111
$ws.Range("H132").Value = 45213.375
$ws.Range("I132").Value = 1019.125
$ws.Range("J132").Value = 177796.12
$ws.Range("K132").Value = 3057.375
$ws.Range("L132").Value = 533388.36
$ws.Range("M132").Value = -527.375
$ws.Range("N132").Value = -538448.36

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 115744.336
$ws.Range("I56").Value = 115744.336
$ws.Range("K56").Value = 115744.336
$ws.Range("M56").Value = -115214.336
$ws.Range("H113").Value = 8575
$ws.Range("I113").Value = 25420.25
$ws.Range("J113").Value = 1088.2222
$ws.Range("K113").Value = 76260.75
$ws.Range("L113").Value = 3264.6666
$ws.Range("M113").Value = -74090.75
$ws.Range("N113").Value = -7604.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 22729736
$ws.Range("I132").Value = 32259562
$ws.Range("K132").Value = 96778686
$ws.Range("M132").Value = -96776156

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2801.3076
$ws.Range("I7").Value = 2037.4546
$ws.Range("K7").Value = 2037.4546
$ws.Range("M7").Value = -1925.4546
$ws.Range("H40").Value = 4846.5
$ws.Range("I40").Value = 2559.4
$ws.Range("K40").Value = 2559.4
$ws.Range("M40").Value = -2423.4
$ws.Range("H93").Value = 1268
$ws.Range("I93").Value = 824.2857
$ws.Range("J93").Value = 1711.7142
$ws.Range("K93").Value = 824.2857
$ws.Range("L93").Value = 1711.7142
$ws.Range("M93").Value = 423.7143
$ws.Range("N93").Value = -4207.7142
$ws.Range("H126").Value = 2801.3076
$ws.Range("I126").Value = 2037.4546
$ws.Range("K126").Value = 6112.3638
$ws.Range("M126").Value = -3642.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1681573
$ws.Range("I122").Value = 2857993.5
$ws.Range("J122").Value = 972.1429000000001
$ws.Range("K122").Value = 8573980.5
$ws.Range("L122").Value = 2916.4287
$ws.Range("M122").Value = -8571530.5
$ws.Range("N122").Value = -7816.4287
$ws.Range("H132").Value = 1145452.2
$ws.Range("I132").Value = 1318472.4
$ws.Range("J132").Value = 3519.6
$ws.Range("K132").Value = 3955417.2
$ws.Range("L132").Value = 10558.8
$ws.Range("M132").Value = -3952887.2
$ws.Range("N132").Value = -15618.8
